$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.526.89'
$ws.Range("E2").Value = '  +0.93%  '

# Row 3
$ws.Range("D3").Value = '1.874.13'
$ws.Range("E3").Value = '  +0.20%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.82'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5083'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.63%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3899'
$ws.Range("D8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08409'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.47%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.104'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.15%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.81'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.35%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.218'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.33%  '

# Row 13
$ws.Range("D13").Value = '1.871.50'
$ws.Range("E13").Value = '  -0.22%  '

# Row 14
$ws.Range("E14").Value = '  +0.42%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.241'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.17%  '

# Row 16
$ws.Range("E16").Value = '  +0.11%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.24%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.24'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.03%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06708'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.17%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.72'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.38%  '

# Row 21
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.933'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.51%  '

# Row 23
$ws.Range("D23").Value = '28.563.09'
$ws.Range("E23").Value = '  +0.91%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.08'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.27%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.233'
$ws.Range("D25").ClearFormats()

# Row 26
$ws.Range("D26").Value = '2.083.30'
$ws.Range("E26").Value = '  -0.18%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.55'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.72%  '

# Row 28
$ws.Range("E28").Value = '  -0.23%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.353'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.08%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.07'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.39%  '

# Row 31
$ws.Range("E31").Value = '  -1.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.043'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.26%  '

# Row 33
$ws.Range("E33").Value = '  -1.47%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.609'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.23%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02453'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06551'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.04%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2160'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.84%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.859'
$ws.Range("D38").ClearFormats()

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.072'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.38%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.251'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.29%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.191'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.64%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6426'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.48%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.08%  '

# Row 44
$ws.Range("E44").Value = '  -0.10%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6039'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.12%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.98'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.42%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.689'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.14%  '

# Row 48
$ws.Range("E48").Value = '  -0.18%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.217'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.31%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.03'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.43%  '

# Row 51
$ws.Range("E51").Value = '  -8.57%  '
